# Weekly update: add two new price records (row 46 and 47) for
# Agricola del Norte S.A. de Arica - Mango, pushing the existing
# data down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 46 (this shifts the old
# rows 46..154 down to 48..156, and keeps row formatting/number
# formats consistent with the surrounding rows).
$ws.Range("A46:A47").EntireRow.Insert()

# --- Row 46 : Especial -------------------------------------------------
$ws.Cells.Item(46,1).Value2  = 1
$ws.Cells.Item(46,2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46,3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(46,4).Value2  = 44883
$ws.Cells.Item(46,5).Value2  = 15
$ws.Cells.Item(46,6).Value2  = "Fruta"
$ws.Cells.Item(46,7).Value2  = 100108
$ws.Cells.Item(46,8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(46,9).Value2  = 100108002
$ws.Cells.Item(46,10).Value2 = "Mango"
$ws.Cells.Item(46,11).Value2 = "Sin especificar"
$ws.Cells.Item(46,12).Value2 = "Especial"
$ws.Cells.Item(46,13).Value2 = 2000
$ws.Cells.Item(46,14).Value2 = 6500
$ws.Cells.Item(46,15).Value2 = 7000
$ws.Cells.Item(46,16).Value2 = 6750
$ws.Cells.Item(46,17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(46,18).Value2 = "Perú"
$ws.Cells.Item(46,19).Value2 = 1688
$ws.Cells.Item(46,20).Value2 = 4

# --- Row 47 : Primera ---------------------------------------------------
$ws.Cells.Item(47,1).Value2  = 1
$ws.Cells.Item(47,2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(47,3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(47,4).Value2  = 44883
$ws.Cells.Item(47,5).Value2  = 15
$ws.Cells.Item(47,6).Value2  = "Fruta"
$ws.Cells.Item(47,7).Value2  = 100108
$ws.Cells.Item(47,8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(47,9).Value2  = 100108002
$ws.Cells.Item(47,10).Value2 = "Mango"
$ws.Cells.Item(47,11).Value2 = "Sin especificar"
$ws.Cells.Item(47,12).Value2 = "Primera"
$ws.Cells.Item(47,13).Value2 = 2800
$ws.Cells.Item(47,14).Value2 = 6500
$ws.Cells.Item(47,15).Value2 = 7000
$ws.Cells.Item(47,16).Value2 = 6679
$ws.Cells.Item(47,17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(47,18).Value2 = "Perú"
$ws.Cells.Item(47,19).Value2 = 1670
$ws.Cells.Item(47,20).Value2 = 4
